$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 20, shifting existing rows 20-28 down to 21-29
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with values (copy of old row 20's data, with
# updated Fecha/price columns reflecting the new weekly entry)
$ws.Range("A20").Value2 = 11
$ws.Range("B20").Value = "Vega Monumental Concepción"
$ws.Range("C20").Value = "Bíobío"
$ws.Range("D20").Value2 = 44755
$ws.Range("E20").Value2 = 8
$ws.Range("F20").Value2 = 100112026
$ws.Range("G20").Value = "Haba"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value2 = 100
$ws.Range("K20").Value2 = 16000
$ws.Range("L20").Value2 = 17000
$ws.Range("M20").Value2 = 16500
$ws.Range("N20").Value = "`$/saco 25 kilos"
$ws.Range("O20").Value = "Región de Coquimbo"
$ws.Range("P20").Value2 = 660
$ws.Range("Q20").Value2 = 25
$ws.Range("R20").Value = "Hortaliza"
